$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1) "Exercises" Heading1 paragraph: add spacing-after=0 and a collapsed
#    "_GoBack" bookmark immediately after the run text (still inside the
#    paragraph, before the paragraph mark).
# -------------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item(1)
$headingPara.Range.ParagraphFormat.SpaceAfter = 0

$endRng = $headingPara.Range.Duplicate
$endRng.End = $endRng.End - 1
$endRng.Collapse(0)
$endRng.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $endRng)
$markerStart = $endRng.Start
$markerEnd = $markerStart + 1
$markerRng = $d.Range($markerStart, $markerEnd)
$markerRng.Delete()

# -------------------------------------------------------------------------
# 2) "Print Stack Usag" + bookmark + "e" -> single run "Print Stack Usage"
#    (the old "_GoBack" bookmark that used to live here is removed as a
#    side effect of replacing the whole matched range).
# -------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Print Stack Usage", $false, $false, $false, $false, $false, $true, 1, $false, "Print Stack Usage", 2)

# -------------------------------------------------------------------------
# 3) Remove the table row for "BLE Low Power (SDS)" (the row right after
#    the "OTA Firmware Upgrade (Secure)" row).
# -------------------------------------------------------------------------
$table = $d.Tables.Item(1)
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows.Item($i)
    $lastCellText = $row.Cells.Item($row.Cells.Count).Range.Text
    if ($lastCellText -like "*BLE Low Power (SDS)*") {
        $row.Delete()
        break
    }
}
